# Rebuild "stadiumArcadium.xlsx" as a proper Album/Song/Lyrics table on a
# sheet renamed to "Sheet2", replacing the old single 12-column, 1-row
# layout. Lyrics text is loaded into single-quoted here-strings (so no
# PowerShell variable/escape expansion touches the song lyrics) and then
# written cell-by-cell, matching the authored workbook's data + wrap-text
# formatting on the Lyrics column.

$s0 = @'
She's only eighteen
Don't like the Rollin' Stones
She took a short cut
To bein' fully grown
She got that mood ring
Her little sister, Rose
A smell of Springsteen
A pair of pantyhose
This talking picture show is leaking from a silhouette
She said, "My man, you know, it's time to get your fingers wet"
You hustle faster when you can't afford a cigarette
The last I heard from you, why, you were screamin', "Handle it!"
Knock the world
Right off its feet and straight onto its head
The book of love will long be laughing after you are dead
Fascinated by the look of you and what was said
Make a play for all the brightest minds in life will share
I heard some P-Funk
Out on the road again
To get your head shrunk is what I recommend
It's in your bloodline
A perfect… 
'@
$s1 = @'
Album
'@
$s2 = @'
Stadium Arcadium
'@
$s3 = @'
Song
'@
$s4 = @'
Dani California
'@
$s5 = @'
Lyrics
'@
$s6 = @'
"Getting born in the state of Mississippi
Papa was a copper, and her mama was a hippy
In Alabama she would swing a hammer
Price you got to pay when you break the panorama
She never knew that there was anything more than poor
What in the world does your company take me for?
Black bandanna, sweet Louisiana
Robbing on a bank in the state of Indiana
She's a runner
Rebel, and a stunner
On her merry way saying baby, watcha gonna?
Looking down the barrel of a hot metal forty-five
Just another way to survive

California, rest in peace
Simultaneous release
California, show your teeth
She's my priestess
I'm your priest
Yeah, yeah, yeah

She's a lover, baby, and a fighter
Should've seen it coming when I got a little brighter
With a name like Dani California
Day was… "
'@
$s7 = @'
Snow
'@
$s8 = @'
"fCome to decide that the things that I tried
 Were in my life just to get high on
 When I sit alone come get a little known
 But I need more than myself this time

Step from the road to the sea to the sky
 And I do believe it, we rely on
 When I lay it on come get to play it on
 All my life to sacrifice

Hey oh, listen what I say, oh
 I got your hey oh, now listen what I say, oh

When will I know that I really can't go
 To the well once more time to decide on
 When it's killing me, when will I really see
 All that I need to look inside
 

Related












Read more: Red Hot Chili Peppers - Snow (hey Oh) Lyrics | MetroLyrics 
"
'@
$s9 = @'
"All aboard, stitch in time
Get yours, got mine
In a minute, I'll be there
Sit tight, get square
You could do it at the hippodrome
Slide back, trombone
Anybody got a TV at home
That's right, unknown

When I pick up on that smell
Pick it up and run like hell
Little woman save me some
Better get up on your run

So much more than 
Charlie's wakin' me
To my core and
Charlie's shakin' me
And tell my story 
And Charlie's makin' me
And Charlie's makin' me smile
Oh oh now

Everybody, do the twist
Get the message, on Flea's fist
Move around like a scientist
Lay down, get kissed
Big picture and it never lies
Big daddy, will advise
Ever ready, in disguise
Sunset, sunrise

When I pick up on that smell
Pick it up and run like hell
Little woman save me some
Better get up on your run

So much… "
'@
$s10 = @'
Charlie
'@
$s11 = @'
Hump De Bump
'@
$s12 = @'
She's Only 18
'@
$s13 = @'
Slow Cheetah
'@
$s14 = @'
Torture Me
'@
$s15 = @'
Strip My Mind
'@
$s16 = @'
Especially in Michigan
'@
$s17 = @'
Warlocks
'@
$s18 = @'
C'mon Girl
'@
$s19 = @'
"
Bells around Saint Petersburg
When I saw you
I hope I get what you deserve
And this is where I find

Smoke surrounds your perfect face
And I'm falling

Pushing a broom out into space
And this where I find the way

The stadium arcadium
A mirror to the moon (a mirror to the moon)
Well I'm forming and I'm warning
State of the art
Until the clouds come crashing

Stranger things have happened
Both before and after noon (before and after noon)
Well I'm forming and I'm warning
Pushin' myself
And no I don't mind asking
Now

Alone inside my forest room
And it's storming

I never thought I'd be in bloom
But this is where I start

Derelict days and the stereo plays
For the all night crowd
That it cannot phase
And I'm calling

Tedious weeds that the media breeds
But the… 






The stadium arcadium
A mirror to the moon (a mirror to the moon)
Well I'm forming and I'm warning
State of the art
Until the clouds come crashing

Stranger things have happened
Both before and after noon (before and after noon)
Well I'm forming and I'm warming
Pushin' myself
And no I don't mind asking
Now

And this is where I find

Rays of dust that wrap around
Your citizen

Kind enough to disavow
And this is where I stand

The stadium arcadium
A mirror to the moon (a mirror to the moon)
Well I'm forming and I'm warning
State of the art
Until the clouds come crashing

Stranger things have happened
Both before and after noon (before and after noon)
And I'm forming and I'm warning
Pushin' myself
And no I don't mind asking

The stadium arcadium
A mirror to the moon (a mirror to the moon)
Well I'm forming and I'm warning
State of the art
Until the clouds come crashing

Stranger things have happened"
'@
$s20 = @'
"Forty detectives this week
Forty detectives strong
Takin' a stroll down love street
Strollin', is that so wrong?
Can I get my co-dependant?

Hump de bump doop bodu
Bump de hump doop bop
Hump de bump doop bodu

Oh no

Bump de hump doop bodu
Hump de bump doop bop
Bump de hump doop bodu

Bump, bump

It must have been a hundred miles
Or any of a hundred styles
It's not about the smile you wear but
The way we make out

When I was an all-aloner
Nothing but a two-beach comber
Anybody seen the sky, I'm
I'm wide awake now!

Workin' the beat as we speak
Workin' the belle du monde
Believe in the havoc we wreak
Believin', is that so wrong?
Can I get my co-dependant?

Hump de bump doop bodu
Bump de hump doop bop
Hump de bump doop bodu

Bump de hump doop bodu
Hump de bump doop bop
Bump de… "
'@
$s21 = @'
"Waking up dead inside of my head
Will never never do there is no med
No medicine to take

I've had a chance to be insane
Asylum from the falling rain
I've had a chance to break

It's so bad it's got to be good
Mysterious girl misunderstood
Dressed like a wedding cake

Any other day and I might play
A funeral march for Bonnie Brae
Why try and run away

Slow cheetah come
Before my forest
Looks like it's on today

Slow cheetah come
It's so euphoric
No matter what they say

I know a girl
She worked in a store
She knew not what
Her life was for
She barely knew her name

They tried to tell her
She would never be
As happy as the girl
In the magazine
She bought it with her pay

Slow cheetah come
Before my forest
Looks like it's on today

Slow cheetah come
It's so euphoric
No matter… "
'@
$s22 = @'
"Because I'm happy to be sad
I want it all I want it bad
Oh oh
It's what I know

A vintage year for pop I hear
The middle of the end is near

Let's go, oh
It's what I know
Torture me and torture me
It's forcin' me so torture me please
Torture me with sorcery
It's forcin' me so torture me please

All the leaves are turning brown
The wind is pushing me around
Let's go
It's what I know

Torture me and torture me
It's forcin' me so torture me please
Torture me with sorcery
It's forcin' me so torture me please

The will of God is standing still
Brazilian children get their fill
Let's go

Let's turn it up and dumb it down
The vision of your ultra sound
Is so

All the leaves are turning brown
The wind is pushing me around
Let's go"
'@
$s23 = @'
"Oh, yeah yeah
Wow, wow, wow, wow, yeah
Arthur J. did, indicate that
The boulevard will never be 
So full of life and love again hey
(Aw, say goodbye to your boots, man)

Hot as Hades, early eighties 
Sing another song
Make me feel like I'm in love again, hey
(You gotta lose to win)

Oh yeah 
Oh 
Please don't strip my mind 
Leave something behind 
Please don't strip my mind 

Oh, hey yeah, oh 
Wow, wow, wow, wow, wow yeah
All in favor sign the waiver 
Bloody Carolina 
Won't you take another look inside, hey
(Aw, it will make me cry)

Operator, co-creator 
Come on baby Ellie 
Won't you blow another compensator, hey
(Ah, you only get what you bring)

Oh yeah 
Oh "
'@
$s24 = @'
"
Life is my friend 
Rake it up to take it in 
Wrap me in your cinnamon 
Especially in Michigan 
Well I could be your friend 

White clouds I'm in 
A mitten full of fisherman 
C'mon huckleberry finn 
Show me how to make her grin 
Well I'm in Michigan 

Cry me a future 
Where the revelations run amok 
Ladies and gentlemen 
Lions and tigers come running 
Just to steal your luck 

A rainy Lithuanian 
Who's dancing as an Indian 
Painted in my tiger skin 
Especially in Michigan

Double chins and bowling pins 
Unholy Presbyterians 
Land is full of medicine 
I find it when I'm slipping in 

The tainted new librarian (into Michigan)
Who fainted when she tucked you in 
Let's float away like zeppelins 
On stoic gusts of northern wind





Cry me a future 
Where the revelations run amok 
Ladies and gentlemen 
Lions and tigers come running 
Just to steal your luck 

Out on the farm we'll be 
Swimming with the mother duck 
Deep in the mitten where 
Lions and tigers come running 
Just to steal your luck 

Life is my friend 
Underwater violins 
Order now from ho chi min 
A porcelain that comes in twins 

Throw me in the looney bin 
'Cause I can take it on the chin (when I'm in Michigan)
The cleavage of your pillow skin 
Is moving like a violin

Cry me a future 
Where the revelations run amok 
Ladies and gentlemen 
Lions and tigers come running 
Just to steal your luck 

Out on the farm we'll be 
Swimming with the mother duck 
Deep in the mitten where 
Lions and tigers come running 
Just to steal your luck, yeah"
'@
$s25 = @'
"
Warlocks in wonderland 
I've gotta megatropolis in my hand 
And a, subterranean marching band 
Makin' noise for the man in the Vatican 
And a 

A little package and off we go 
Oh, ticky ticky tackita tic tac toe 
I know, everybody's Eskimo 
We've got another thing coming 
And that's our show, well 

Every night I go looking for you 
Everyone in the world adores you 
A little pocket of something kind 
To find your reason 
Coming up on it everyday for 
Look at me and it's what I stay for 
A little locket of fantasy 
That we believe in 

Lilacs and contraband 
I've got Santa Monica in my hand 
A little, Beatlemania when I can 
And I've got two big bags of old Japan 






Ring side and blow-by-blow 
Another, main event at the old rainbow 
We're comin', right on top of the tupelo 
When she looks just like Brigitte Bardot 

Every night I go looking for you 
Everyone in the world adores you 
A little pocket of something kind 
To find your reason 
Coming up on it every day for 
Look at me and it's what I stay for 
A little locket of fantasy 
That we believe in 

Make a deal with Uncle Weezer 
Sign your name to claim 
China Chow will try to please her 
Sweetness came from Jane 

Warlocks in wonderland 
I've gotta, Rockapotamus in my hand 
With a, happy ending that's made of sand 
With a little bit of lovin' is all I can

Every night I go looking for you 
Everyone in the world adores you 
A little pocket of something kind 
To find your reason 
Coming up on it every day for 
Look at me and it's what I stay for 
A little locket of fantasy 
That we believe in 

Every night I go looking for you 
Everyone in the world adores you 
A little pocket of something kind 
To find your reason 
Coming up on it every day for 
Look at me and it's what I stay for 
A little locket of fantasy 
That we believe in"
'@
$s26 = @'
"
Everything inside of me
Is burning up for you to see
And if we should get into it
A two-way kind of syndicate
Hold me

Compare us to a Ferris wheel
Just be sure what you want to steal
Some people do it by the book
But I prefer to go by feel
Tell me

I wanted to get out of here
But every time I reappear
Now I have the words my dear
To whisper out into your ear
Hey, yeah

C'mon girl (let's go)
Let's get it right (let's go)
Let's get it right (let's go)
Let's get it right (let's go)

C'mon girl (let's go)
Let's get it right (let's go)
Let's get it right (let's go)
Let's get it right (let's go)

C'mon girl (let's go)
Let's get it right (let's go)
Let's get it right (let's go)
Let's get it right (let's go)






C'mon girl (let's go)
Let's get it right (let's go)
Let's get it right (let's go)
Let's get it right (let's go)

The spirit of a dragon's tear
Is lovely at this time of year
The cave within your mountainside
Is deeper than it will be wide
Hold me

My disenchanted diplomat
Asleep inside the laundromat
Conveyor belts are moving
And I want you to be sure of that
Tell me

And if you let it germinate
I know it will be worth the wait
Disputed by the news that it was
You who tried to instigate
Hey, yeah

C'mon girl (let's go)
Let's get it right (let's go)
Let's get it right (let's go)
Let's get it right (let's go)

C'mon girl (let's go)
Let's get it right (let's go)
Let's get it right (let's go)
Let's get it right (let's go)

C'mon girl (let's go)
Let's get it right (let's go)
Let's get it right (let's go)
Let's get it right (let's go)

C'mon girl (let's go)
Let's get it right (let's go)
Let's get it right (let's go)
Let's get it right (let's go)

Blessed are the hypocrites
Outwitted but she never quits
The trouble with the band of slits
Is washing off the muddy bits
Hey, yeah

C'mon girl (let's go)
Let's get it right (let's go)
Let's get it right (let's go)
Let's get it right (let's go)

C'mon girl (let's go)
Let's get it right (let's go)
Let's get it right (let's go)
Let's get it right (let's go)

C'mon girl (let's go)
Let's get it right (let's go)
Let's get it right (let's go)
Let's get it right (let's go)

C'mon girl (let's go)
Let's get it right (let's go)
Let's get it right (let's go)
Let's get it right (let's go)

Oh yeah she's with me and I'm your man
If I can't find you no one can"
'@

# Start from a clean slate: remove all existing rows/content & formatting
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Sheet2"
if ($ws.UsedRange -ne $null) {
    $ws.UsedRange.EntireRow.Delete()
}

# Header row
$ws.Range("A1").Value = $s1
$ws.Range("B1").Value = $s3
$ws.Range("C1").Value = $s5

# Data rows 2..13
$ws.Range("A2").Value = $s2
$ws.Range("B2").Value = $s4
$ws.Range("C2").Value = $s6

$ws.Range("A3").Value = $s2
$ws.Range("B3").Value = $s7
$ws.Range("C3").Value = $s8

$ws.Range("A4").Value = $s2
$ws.Range("B4").Value = $s10
$ws.Range("C4").Value = $s9

$ws.Range("A5").Value = $s2
$ws.Range("B5").Value = $s2
$ws.Range("C5").Value = $s19

$ws.Range("A6").Value = $s2
$ws.Range("B6").Value = $s11
$ws.Range("C6").Value = $s20

$ws.Range("A7").Value = $s2
$ws.Range("B7").Value = $s12
$ws.Range("C7").Value = $s0

$ws.Range("A8").Value = $s2
$ws.Range("B8").Value = $s13
$ws.Range("C8").Value = $s21

$ws.Range("A9").Value = $s2
$ws.Range("B9").Value = $s14
$ws.Range("C9").Value = $s22

$ws.Range("A10").Value = $s2
$ws.Range("B10").Value = $s15
$ws.Range("C10").Value = $s23

$ws.Range("A11").Value = $s2
$ws.Range("B11").Value = $s16
$ws.Range("C11").Value = $s24

$ws.Range("A12").Value = $s2
$ws.Range("B12").Value = $s17
$ws.Range("C12").Value = $s25

$ws.Range("A13").Value = $s2
$ws.Range("B13").Value = $s18
$ws.Range("C13").Value = $s26

# Apply wrap text + row height to the Lyrics column (C2:C13), matching the
# original workbook's style (wrapText, very tall rows)
$lyricsRange = $ws.Range("C2:C13")
$lyricsRange.WrapText = $true
for ($r = 2; $r -le 13; $r++) {
    $ws.Rows.Item($r).RowHeight = 409
}

# Selection matches the authored file (active cell B13)
$ws.Range("B13").Select()

Write-Output "Table build complete"
